# Auto-generated Excel COM-interop script applying the Typhon_Profits.xlsx diff
# (numeric recalculated-profit figures across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1516695.1
$ws.Range("J17").Value = 1564070.6
$ws.Range("L17").Value = 4692211.800000001
$ws.Range("N17").Value = -4692547.800000001
$ws.Range("H42").Value = 117.3
$ws.Range("J42").Value = 119.22222
$ws.Range("L42").Value = 357.66666
$ws.Range("N42").Value = -817.66666
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 6000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -6384
$ws.Range("H113").Value = 47623012
$ws.Range("I113").Value = 76926340
$ws.Range("J113").Value = 5108.5
$ws.Range("K113").Value = 76926340
$ws.Range("L113").Value = 5108.5
$ws.Range("M113").Value = -76923086
$ws.Range("N113").Value = -11616.5
$ws.Range("H125").Value = 1066.6666
$ws.Range("I125").Value = 1100
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 9900
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -7440
$ws.Range("N125").Value = -13920
$ws.Range("H129").Value = 196891.28
$ws.Range("J129").Value = 213616.06
$ws.Range("L129").Value = 640848.1799999999
$ws.Range("N129").Value = -650848.1799999999
$ws.Range("H132").Value = 7383.3335
$ws.Range("I132").Value = 10655.429
$ws.Range("J132").Value = 2802.4
$ws.Range("K132").Value = 31966.287
$ws.Range("L132").Value = 8407.200000000001
$ws.Range("M132").Value = -29436.287
$ws.Range("N132").Value = -13467.2
$ws.Range("H138").Value = 1358.6
$ws.Range("I138").Value = 557
$ws.Range("J138").Value = 1988.4286
$ws.Range("K138").Value = 1671
$ws.Range("L138").Value = 5965.2858
$ws.Range("M138").Value = 3469
$ws.Range("N138").Value = -16245.2858
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1455.5
$ws.Range("I2").Value = 1261
$ws.Range("J2").Value = 2168.6667
$ws.Range("K2").Value = 1261
$ws.Range("L2").Value = 2168.6667
$ws.Range("M2").Value = -1148
$ws.Range("N2").Value = -2394.6667
$ws.Range("H32").Value = 29047.795
$ws.Range("I32").Value = 28890.79
$ws.Range("K32").Value = 28890.79
$ws.Range("M32").Value = -28603.79
$ws.Range("H61").Value = 2263
$ws.Range("I61").Value = 1678.28
$ws.Range("J61").Value = 3591.9092
$ws.Range("K61").Value = 1678.28
$ws.Range("L61").Value = 3591.9092
$ws.Range("M61").Value = -1466.28
$ws.Range("N61").Value = -4015.9092
$ws.Range("H116").Value = 1455.5
$ws.Range("I116").Value = 1261
$ws.Range("J116").Value = 2168.6667
$ws.Range("K116").Value = 1261
$ws.Range("L116").Value = 2168.6667
$ws.Range("M116").Value = 1033
$ws.Range("N116").Value = -6756.6667
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 14738.975
$ws.Range("I132").Value = 1957.8214
$ws.Range("J132").Value = 44561.668
$ws.Range("K132").Value = 5873.4642
$ws.Range("L132").Value = 133685.004
$ws.Range("M132").Value = -3343.4642
$ws.Range("N132").Value = -138745.004
$ws.Range("H136").Value = 2263
$ws.Range("I136").Value = 1678.28
$ws.Range("J136").Value = 3591.9092
$ws.Range("K136").Value = 5034.84
$ws.Range("L136").Value = 10775.7276
$ws.Range("M136").Value = -2484.84
$ws.Range("N136").Value = -15875.7276
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1455.5
$ws.Range("I3").Value = 1261
$ws.Range("J3").Value = 2168.6667
$ws.Range("K3").Value = 1261
$ws.Range("L3").Value = 2168.6667
$ws.Range("M3").Value = -1147
$ws.Range("N3").Value = -2396.6667
$ws.Range("H99").Value = 1186.85
$ws.Range("I99").Value = 827.3125
$ws.Range("K99").Value = 827.3125
$ws.Range("M99").Value = 670.6875
$ws.Range("H134").Value = 26467.373
$ws.Range("I134").Value = 32972.94
$ws.Range("K134").Value = 98918.82000000001
$ws.Range("M134").Value = -96383.82000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 86142856
$ws.Range("I6").Value = 12599998
$ws.Range("K6").Value = 12599998
$ws.Range("M6").Value = -12599885
$ws.Range("H21").Value = 14571.429
$ws.Range("I21").Value = 9000
$ws.Range("J21").Value = 15500
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 15500
$ws.Range("M21").Value = -8765
$ws.Range("N21").Value = -15970
$ws.Range("H31").Value = 14785.971
$ws.Range("I31").Value = 30165.143
$ws.Range("J31").Value = 4020.55
$ws.Range("K31").Value = 30165.143
$ws.Range("L31").Value = 4020.55
$ws.Range("M31").Value = -29870.143
$ws.Range("N31").Value = -4610.55
$ws.Range("H34").Value = 14785.971
$ws.Range("I34").Value = 30165.143
$ws.Range("J34").Value = 4020.55
$ws.Range("K34").Value = 30165.143
$ws.Range("L34").Value = 4020.55
$ws.Range("M34").Value = -29963.143
$ws.Range("N34").Value = -4424.55
$ws.Range("H58").Value = 13507.725
$ws.Range("I58").Value = 1069.1852
$ws.Range("J58").Value = 39341.617
$ws.Range("K58").Value = 1069.1852
$ws.Range("L58").Value = 39341.617
$ws.Range("M58").Value = -866.1851999999999
$ws.Range("N58").Value = -39747.617
$ws.Range("H99").Value = 17861194
$ws.Range("I99").Value = 3339
$ws.Range("J99").Value = 41671668
$ws.Range("K99").Value = 3339
$ws.Range("L99").Value = 41671668
$ws.Range("M99").Value = -1841
$ws.Range("N99").Value = -41674664
$ws.Range("H103").Value = 25577
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 25577
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 25577
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -27921
$ws.Range("H126").Value = 17861194
$ws.Range("I126").Value = 3339
$ws.Range("J126").Value = 41671668
$ws.Range("K126").Value = 10017
$ws.Range("L126").Value = 125015004
$ws.Range("M126").Value = -7547
$ws.Range("N126").Value = -125019944
$ws.Range("H134").Value = 1140.7255
$ws.Range("I134").Value = 862.2083
$ws.Range("J134").Value = 1388.2963
$ws.Range("K134").Value = 2586.6249
$ws.Range("L134").Value = 4164.8889
$ws.Range("M134").Value = -51.6248999999998
$ws.Range("N134").Value = -9234.8889
$ws.Range("H136").Value = 13507.725
$ws.Range("I136").Value = 1069.1852
$ws.Range("J136").Value = 39341.617
$ws.Range("K136").Value = 3207.5556
$ws.Range("L136").Value = 118024.851
$ws.Range("M136").Value = -657.5555999999997
$ws.Range("N136").Value = -123124.851
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1286.7778
$ws.Range("J5").Value = 1567.75
$ws.Range("L5").Value = 4703.25
$ws.Range("N5").Value = -4927.25
$ws.Range("H36").Value = 151901.17
$ws.Range("J36").Value = 225950.75
$ws.Range("L36").Value = 677852.25
$ws.Range("N36").Value = -678190.25
$ws.Range("H37").Value = 22799546
$ws.Range("J37").Value = 22799546
$ws.Range("L37").Value = 68398638
$ws.Range("N37").Value = -68398862
$ws.Range("H86").Value = 100000376
$ws.Range("J86").Value = 125000216
$ws.Range("L86").Value = 375000648
$ws.Range("N86").Value = -375003020
$ws.Range("H89").Value = 100000376
$ws.Range("J89").Value = 125000216
$ws.Range("L89").Value = 1125001944
$ws.Range("N89").Value = -1125013800
$ws.Range("H107").Value = 8594.833000000001
$ws.Range("I107").Value = 20220
$ws.Range("K107").Value = 60660
$ws.Range("M107").Value = -58740
$ws.Range("H122").Value = 697.9375
$ws.Range("I122").Value = 377.14285
$ws.Range("J122").Value = 947.44446
$ws.Range("K122").Value = 3394.28565
$ws.Range("L122").Value = 8527.00014
$ws.Range("M122").Value = -944.2856500000003
$ws.Range("N122").Value = -13427.00014
$ws.Range("H131").Value = 788.0404
$ws.Range("J131").Value = 792.875
$ws.Range("L131").Value = 2378.625
$ws.Range("N131").Value = -12458.625
$ws.Range("H132").Value = 1451.9
$ws.Range("I132").Value = 1404.8334
$ws.Range("K132").Value = 12643.5006
$ws.Range("M132").Value = -10113.5006
$ws.Range("H135").Value = 1286.7778
$ws.Range("J135").Value = 1567.75
$ws.Range("L135").Value = 14109.75
$ws.Range("N135").Value = -19179.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 21741144
$ws.Range("I102").Value = 22729342
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 22729342
$ws.Range("L102").Value = 800
$ws.Range("M102").Value = -22727720
$ws.Range("N102").Value = -4044
$ws.Range("H123").Value = 7351.7393
$ws.Range("I123").Value = 3573.3333
$ws.Range("J123").Value = 9780.714
$ws.Range("K123").Value = 3573.3333
$ws.Range("L123").Value = 9780.714
$ws.Range("M123").Value = -1123.3333
$ws.Range("N123").Value = -14680.714
$ws.Range("H126").Value = 4213.5454
$ws.Range("I126").Value = 3396.4736
$ws.Range("J126").Value = 5322.4287
$ws.Range("K126").Value = 10189.4208
$ws.Range("L126").Value = 15967.2861
$ws.Range("M126").Value = -7719.4208
$ws.Range("N126").Value = -20907.2861
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 499.27274
$ws.Range("I16").Value = 499
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 499
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -329
$ws.Range("N16").Value = -840
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 6333.3335
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 13000
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 13000
$ws.Range("N33").Value = -13500
$ws.Range("M33").Value = -2750
$ws.Range("H36").Value = 6333.3335
$ws.Range("I36").Value = 3000
$ws.Range("J36").Value = 13000
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 13000
$ws.Range("N36").Value = -13500
$ws.Range("M36").Value = -2750
$ws.Range("H62").Value = 3645.375
$ws.Range("J62").Value = 4156.3335
$ws.Range("L62").Value = 4156.3335
$ws.Range("N62").Value = -5404.3335
$ws.Range("H65").Value = 3645.375
$ws.Range("J65").Value = 4156.3335
$ws.Range("L65").Value = 20781.6675
$ws.Range("N65").Value = -27021.6675

Write-Host "Applied $(262) value updates and $(3) clears."